# Insert 3 new rows before row 1060, shifting the existing rows 1060-1084
# down to 1063-1087, then populate the 3 newly inserted rows with the new
# weekly price records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at position 1060 (pushes old 1060..1084 -> 1063..1087)
$ws.Rows("1060:1062").Insert()

# --- New row 1060 : Cuatro cascos verde ---
$ws.Cells.Item(1060, 1).Value = 5
$ws.Cells.Item(1060, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(1060, 3).Value = "Maule"
$ws.Cells.Item(1060, 4).Value = 45239
$ws.Cells.Item(1060, 5).Value = 7
$ws.Cells.Item(1060, 6).Value = 100112002
$ws.Cells.Item(1060, 7).Value = "Pimiento"
$ws.Cells.Item(1060, 8).Value = "Cuatro cascos verde"
$ws.Cells.Item(1060, 9).Value = "Primera"
$ws.Cells.Item(1060, 10).Value = 100
$ws.Cells.Item(1060, 11).Value = 28000
$ws.Cells.Item(1060, 12).Value = 28000
$ws.Cells.Item(1060, 13).Value = 28000
$ws.Cells.Item(1060, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(1060, 15).Value = "Regi" + [char]0x00F3 + "n del Maule"
$ws.Cells.Item(1060, 16).Value = 1867
$ws.Cells.Item(1060, 17).Value = 15
$ws.Cells.Item(1060, 18).Value = "Hortaliza"

# --- New row 1061 : Zafiro rojo ---
$ws.Cells.Item(1061, 1).Value = 5
$ws.Cells.Item(1061, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(1061, 3).Value = "Maule"
$ws.Cells.Item(1061, 4).Value = 45239
$ws.Cells.Item(1061, 5).Value = 7
$ws.Cells.Item(1061, 6).Value = 100112002
$ws.Cells.Item(1061, 7).Value = "Pimiento"
$ws.Cells.Item(1061, 8).Value = "Zafiro rojo"
$ws.Cells.Item(1061, 9).Value = "Primera"
$ws.Cells.Item(1061, 10).Value = 200
$ws.Cells.Item(1061, 11).Value = 30000
$ws.Cells.Item(1061, 12).Value = 30000
$ws.Cells.Item(1061, 13).Value = 30000
$ws.Cells.Item(1061, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(1061, 15).Value = "Regi" + [char]0x00F3 + "n de Arica y Parinacota"
$ws.Cells.Item(1061, 16).Value = 2000
$ws.Cells.Item(1061, 17).Value = 15
$ws.Cells.Item(1061, 18).Value = "Hortaliza"

# --- New row 1062 : Zafiro verde ---
$ws.Cells.Item(1062, 1).Value = 5
$ws.Cells.Item(1062, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(1062, 3).Value = "Maule"
$ws.Cells.Item(1062, 4).Value = 45239
$ws.Cells.Item(1062, 5).Value = 7
$ws.Cells.Item(1062, 6).Value = 100112002
$ws.Cells.Item(1062, 7).Value = "Pimiento"
$ws.Cells.Item(1062, 8).Value = "Zafiro verde"
$ws.Cells.Item(1062, 9).Value = "Primera"
$ws.Cells.Item(1062, 10).Value = 200
$ws.Cells.Item(1062, 11).Value = 30000
$ws.Cells.Item(1062, 12).Value = 30000
$ws.Cells.Item(1062, 13).Value = 30000
$ws.Cells.Item(1062, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(1062, 15).Value = "Regi" + [char]0x00F3 + "n de Arica y Parinacota"
$ws.Cells.Item(1062, 16).Value = 2000
$ws.Cells.Item(1062, 17).Value = 15
$ws.Cells.Item(1062, 18).Value = "Hortaliza"
